# Applies the diff:
#  1. Append three new paragraphs to the Content Placeholder of slide 11
#     ("Modelling Digital Twin Data"), right after the existing
#     "... but fall short in capturing the complex inter-entity dynamics."
#     paragraph.
#  2. Delete slide 12 (the slide that only held those three paragraphs:
#     "Even our Data Platform Design methodology ...", "Yet, no multi-store
#     solution ...", "What about an hybrid data structure?") now that its
#     text lives on slide 11 instead.

$p = $ppt.ActivePresentation

# --- Step 1: extend slide 11's content placeholder text -------------------
$s11 = $p.Slides.Item(11)
$contentShape = $s11.Shapes.Item(2)

$tr = $contentShape.TextFrame.TextRange
$beforeLen = $tr.Length

$para1 = "Even the Data Platform Design methodology suggested different architectures tailored to each DT…"
$para2 = "Yet, no multi-store solution has yet achieved broad adoption in the literature."
$para3 = "What about an hybrid data structure?"

[void]$tr.InsertAfter("`r" + $para1 + "`r" + $para2 + "`r" + $para3)

$full = $contentShape.TextFrame.TextRange

$start1 = $beforeLen + 2
$len1 = $para1.Length
$start2 = $start1 + $len1 + 1
$len2 = $para2.Length
$start3 = $start2 + $len2 + 1
$len3 = $para3.Length

$r1 = $full.Characters($start1, $len1)
$r1.IndentLevel = 0

$r2 = $full.Characters($start2, $len2)
$r2.IndentLevel = 0

$r3 = $full.Characters($start3, $len3)
$r3.IndentLevel = 0
$r3.Font.Bold = $true

# --- Step 2: remove the now-redundant slide that used to carry this text --
[void]$p.Slides.Item(12).Delete()
